$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Record", "RJ Record", "Governo", "2025-04-07T18:11", "Neutro", "Ato pacífico. Movimento dos Sem Terra faz protesto no acesso à usina Sapucaia, em Campos. Repórter *ao vivo*. Entrevista com diretora nacional do MST, Ero Silva, e com comandante da GCM, Wellington Levino. "),
    @("Record", "RJ Record", "Economia", "2025-04-07T18:21", "Positivo", "Campos vai realizar o Mutirão ""Facilita Imposto de Renda"". Repórter *ao vivo*. "),
    @("Record", "RJ Record", "Saúde", "2025-04-07T18:22", "Positivo", "Gripe Influenza. Campos inicia campanha de imunização contra o vírus.  Entrevista com dois idosos que foram se vacinar e com o subsecretário de Vigilância em Saúde, Charbell Kury. Meta é imunizar 90% do público alvo. *matéria*"),
    @("Record", "RJ Record", "Defesa Civil", "2025-04-07T18:33", "Neutro", "Efeito frente fria. Alerta da Marinha se confirma e ressaca no mar atinge a região. Repórter *ao vivo*. Em Farol, as ondas atingiram 3 m de altura. Pescadores suspenderam as atividades. "),
    @("Record", "RJ Record", "Obras", "2025-04-07T19:08", "Neutro", "Parte dos prédios históricos da cidade ainda precisam de manutenção. Entrevista com chefe do Inepac NF, Geovani Laurindo, e com a historiadora Sylvia Paes. Necessidade de 434 prédios tombados em Campos, sendo a maioria na área central. Semana passada, Defesa Civil começou a demolir prédio que estava ameaçado. Proprietários e poder público precisam cuidar. Representante do Inepac falou que Campos oferece desconto do IPTU. *matéria* também foi exibida no Balanço Geral. ")
)

$startRow = 70
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
